# The weekly refresh reshuffles which sample row carries which
# Mercado/Fecha/Calidad/Volumen/Precio data. The Fecha (D), Calidad (I),
# Volumen (J), Precio minimo (K), Precio maximo (L), Precio promedio
# ponderado (M) and Precio $/Kg (P) columns are permuted across rows
# while the remaining columns (Mercado, Region, Categoria, Variedad,
# Unidad de comercializacion, Origen, Kg o Unidades, Clasificacion) stay
# put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (values currently sitting in the source
# row move into the destination row)
$rowMap = @{
    2  = 7
    3  = 11
    6  = 8
    7  = 9
    8  = 6
    9  = 12
    10 = 18
    11 = 17
    12 = 19
    13 = 20
    16 = 10
    17 = 2
    18 = 3
    19 = 16
    20 = 13
}

$cols = @("D", "I", "J", "K", "L", "M", "P")

# Snapshot every affected cell's current value before writing anything,
# so overlapping source/destination rows never read back an already
# overwritten value.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Range("$col$srcRow").Value2
        }
        $snapshot[$srcRow] = $rowVals
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowVals[$col]
    }
}
